# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 4
# (the e56b9177-... entry) for both the zh-cn and de-de sheets,
# as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-12 08:30:54"
$wsZhCn.Range("H4").Value = "2016-03-12 08:31:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-12 08:30:58"
$wsDeDe.Range("H4").Value = "2016-03-12 08:31:29"
